$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9478
$ws.Range("C2").Value = 9456
$ws.Range("D2").Value = 8369
$ws.Range("E2").Value = 0.8850465313028765
$ws.Range("F2").Value = 0.8829921924456636
$ws.Range("G2").Value = 0.09597596898787557
$ws.Range("H2").Value = 0.08474603127870128
$ws.Range("I2").Value = 41011357.62888187
$ws.Range("J2").Value = 14306014.57377693
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 14306014.57377693
$ws.Range("M2").Value = 55317372.20265881
$ws.Range("N2").Value = 812856555.5672001
$ws.Range("O2").Value = 793677941.7132001
$ws.Range("P2").Value = 0.01759967915100887
$ws.Range("Q2").Value = 0.01802496179104659

# Row 3
$ws.Range("B3").Value = 9762
$ws.Range("C3").Value = 9743
$ws.Range("D3").Value = 8661
$ws.Range("E3").Value = 0.8889459098840193
$ws.Range("F3").Value = 0.8872157344806392
$ws.Range("G3").Value = 0.09352092230585833
$ws.Range("H3").Value = 0.08297323377289889
$ws.Range("I3").Value = 43255456.65684056
$ws.Range("J3").Value = 15158002.11521988
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 15158002.11521988
$ws.Range("M3").Value = 58413458.77206045
$ws.Range("N3").Value = 868756639.5402131
$ws.Range("O3").Value = 848873415.5163431
$ws.Range("P3").Value = 0.01744792664058626
$ws.Range("Q3").Value = 0.01785661070090143

# Row 4
$ws.Range("B4").Value = 10048
$ws.Range("C4").Value = 10026
$ws.Range("D4").Value = 8905
$ws.Range("E4").Value = 0.8881907041691602
$ws.Range("F4").Value = 0.8862460191082803
$ws.Range("G4").Value = 0.09177983476932886
$ws.Range("H4").Value = 0.08133951319873342
$ws.Range("I4").Value = 45461625.99696768
$ws.Range("J4").Value = 15941015.34003143
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 15941015.34003143
$ws.Range("M4").Value = 61402641.33699911
$ws.Range("N4").Value = 923597016.4015658
$ws.Range("O4").Value = 902839601.727118
$ws.Range("P4").Value = 0.01725970857088664
$ws.Range("Q4").Value = 0.01765653091594179

# Row 5
$ws.Range("B5").Value = 10340
$ws.Range("C5").Value = 10303
$ws.Range("D5").Value = 9132
$ws.Range("E5").Value = 0.8863437833640687
$ws.Range("F5").Value = 0.8831721470019342
$ws.Range("G5").Value = 0.09026420469662112
$ws.Range("H5").Value = 0.07971883145933695
$ws.Range("I5").Value = 47924844.26765846
$ws.Range("J5").Value = 16839396.97316202
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 16839396.97316202
$ws.Range("M5").Value = 64764241.24082048
$ws.Range("N5").Value = 981614581.1341684
$ws.Range("O5").Value = 958614956.0512464
$ws.Range("P5").Value = 0.01715479506600808
$ws.Range("Q5").Value = 0.01756638248429519

# Row 6
$ws.Range("B6").Value = 10667
$ws.Range("C6").Value = 10640
$ws.Range("D6").Value = 9432
$ws.Range("E6").Value = 0.8864661654135338
$ws.Range("F6").Value = 0.8842223680509984
$ws.Range("G6").Value = 0.08827800578727582
$ws.Range("H6").Value = 0.07805738732404476
$ws.Range("I6").Value = 50668659.16144493
$ws.Range("J6").Value = 17839816.11902389
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 17839816.11902389
$ws.Range("M6").Value = 68508475.28046882
$ws.Range("N6").Value = 1049528624.192984
$ws.Range("O6").Value = 1025805590.179856
$ws.Range("P6").Value = 0.01699793193610273
$ws.Range("Q6").Value = 0.01739103031783636

Write-Host "Update complete"